$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.294.98"
$ws.Cells.Item(2, 5).Value = "  +0.51%  "
$ws.Cells.Item(3, 4).Value = "1.775.43"
$ws.Cells.Item(3, 5).Value = "  +3.44%  "
$ws.Cells.Item(5, 4).Value = "'313.85"
$ws.Cells.Item(5, 5).Value = "  +0.81%  "
$ws.Cells.Item(6, 4).Value = "'1.001"
$ws.Cells.Item(6, 5).Value = "  +0.06%  "
$ws.Cells.Item(7, 4).Value = "'0.5152"
$ws.Cells.Item(7, 5).Value = "  +7.16%  "
$ws.Cells.Item(8, 4).Value = "'0.3687"
$ws.Cells.Item(8, 5).Value = "  +6.59%  "
$ws.Cells.Item(9, 4).Value = "'42.64"
$ws.Cells.Item(9, 5).Value = "  -0.27%  "
$ws.Cells.Item(10, 4).Value = "'0.07396"
$ws.Cells.Item(10, 5).Value = "  +1.34%  "
$ws.Cells.Item(11, 5).Value = "  +3.90%  "
$ws.Cells.Item(12, 5).Value = "  +0.10%  "
$ws.Cells.Item(13, 5).Value = "  +2.73%  "
$ws.Cells.Item(14, 4).Value = "'6.068"
$ws.Cells.Item(14, 5).Value = "  +3.23%  "
$ws.Cells.Item(15, 4).Value = "1.768.89"
$ws.Cells.Item(15, 5).Value = "  +3.08%  "
$ws.Cells.Item(16, 4).Value = "'6.962"
$ws.Cells.Item(16, 5).Value = "  +1.18%  "
$ws.Cells.Item(17, 4).Value = "'89.23"
$ws.Cells.Item(17, 5).Value = "  +0.01%  "
$ws.Cells.Item(18, 5).Value = "  +0.43%  "
$ws.Cells.Item(19, 4).Value = "'0.06431"
$ws.Cells.Item(19, 5).Value = "  +1.03%  "
$ws.Cells.Item(20, 4).Value = "'1.000"
$ws.Cells.Item(20, 5).Value = "  +0.05%  "
$ws.Cells.Item(21, 5).Value = "  +1.33%  "
$ws.Cells.Item(22, 4).Value = "'5.825"
$ws.Cells.Item(22, 5).Value = "  +3.06%  "
$ws.Cells.Item(23, 4).Value = "27.338.77"
$ws.Cells.Item(23, 5).Value = "  +0.48%  "
$ws.Cells.Item(24, 4).Value = "'11.24"
$ws.Cells.Item(24, 5).Value = "  +3.61%  "
$ws.Cells.Item(25, 4).Value = "'2.117"
$ws.Cells.Item(25, 5).Value = "  +1.37%  "
$ws.Cells.Item(26, 4).Value = "'154.23"
$ws.Cells.Item(26, 5).Value = "  +1.25%  "
$ws.Cells.Item(27, 4).Value = "'20.21"
$ws.Cells.Item(27, 5).Value = "  +2.51%  "
$ws.Cells.Item(28, 4).Value = "'2.335"
$ws.Cells.Item(28, 5).Value = "  +10.91%  "
$ws.Cells.Item(29, 4).Value = "1.976.46"
$ws.Cells.Item(29, 5).Value = "  +3.41%  "
$ws.Cells.Item(30, 4).Value = "'121.15"
$ws.Cells.Item(30, 5).Value = "  +0.61%  "
$ws.Cells.Item(31, 4).Value = "'1.065"
$ws.Cells.Item(31, 5).Value = "  +3.98%  "
$ws.Cells.Item(32, 4).Value = "'0.09781"
$ws.Cells.Item(32, 5).Value = "  +5.60%  "
$ws.Cells.Item(33, 4).Value = "'5.576"
$ws.Cells.Item(33, 5).Value = "  +4.08%  "
$ws.Cells.Item(34, 5).Value = "  +1.31%  "
$ws.Cells.Item(35, 4).Value = "'0.02249"
$ws.Cells.Item(35, 5).Value = "  +1.99%  "
$ws.Cells.Item(36, 5).Value = "  +0.61%  "
$ws.Cells.Item(37, 4).Value = "'11.28"
$ws.Cells.Item(37, 5).Value = "  +1.24%  "
$ws.Cells.Item(38, 4).Value = "'0.6154"
$ws.Cells.Item(38, 5).Value = "  +3.34%  "
$ws.Cells.Item(39, 4).Value = "'4.840"
$ws.Cells.Item(39, 5).Value = "  +1.42%  "
$ws.Cells.Item(40, 4).Value = "'0.2022"
$ws.Cells.Item(40, 5).Value = "  +0.90%  "
$ws.Cells.Item(41, 5).Value = "  +0.99%  "
$ws.Cells.Item(42, 4).Value = "'8.102"
$ws.Cells.Item(42, 5).Value = "  +8.15%  "
$ws.Cells.Item(43, 4).Value = "'1.136"
$ws.Cells.Item(43, 5).Value = "  +3.43%  "
$ws.Cells.Item(44, 4).Value = "'13.09"
$ws.Cells.Item(44, 5).Value = "  +3.52%  "
$ws.Cells.Item(45, 4).Value = "'0.5774"
$ws.Cells.Item(45, 5).Value = "  +2.48%  "
$ws.Cells.Item(46, 4).Value = "'3.636"
$ws.Cells.Item(46, 5).Value = "  +1.21%  "
$ws.Cells.Item(47, 4).Value = "'121.73"
$ws.Cells.Item(47, 5).Value = "  +2.52%  "
$ws.Cells.Item(48, 4).Value = "'1.890"
$ws.Cells.Item(49, 4).Value = "'1.115"
$ws.Cells.Item(49, 5).Value = "  +2.62%  "
$ws.Cells.Item(50, 4).Value = "'0.06718"
$ws.Cells.Item(51, 4).Value = "'70.67"
$ws.Cells.Item(51, 5).Value = "  +1.19%  "
